$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.685.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.072.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "58.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.377.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.169.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.626.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -5.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0977"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.454.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.64%  "
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.263.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.69%  "
